# ATSG quarterly financials update: add the two newest quarters as new
# columns D and E, shifting the existing quarterly history two columns to the
# right (old D:K -> F:M), and correct a handful of restated prior-quarter
# figures that now live in columns H:I.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns before column D; this pushes the existing quarter
# columns (old D:K) two places over to F:M, matching the new A5:M102 extent.
$ws.Columns("D:E").Insert()

# The new blank D:E columns default to General formatting; copy the number /
# date formatting from column F (the old column D, now shifted one quarter
# over) into D:E for each of the three data blocks on the sheet so the new
# cells pick up the same styles as the rest of the grid.
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new quarter columns: D = newest quarter, E = prior quarter.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 280800
$ws.Range("E8").Value = 204900
$ws.Range("D9").Value = 108300
$ws.Range("E9").Value = 78500
$ws.Range("D10").Value = 172600
$ws.Range("E10").Value = 126400
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 54100
$ws.Range("E15").Value = 43200
$ws.Range("D17").Value = 248100
$ws.Range("E17").Value = 178100
$ws.Range("D18").Value = 32700
$ws.Range("E18").Value = 26800
$ws.Range("D20").Value = -22100
$ws.Range("E20").Value = 17400
$ws.Range("D21").Value = 71600
$ws.Range("E21").Value = 93400
$ws.Range("D22").Value = 12500
$ws.Range("E22").Value = 5600
$ws.Range("D23").Value = -1900
$ws.Range("E23").Value = 38600
$ws.Range("D24").Value = 3300
$ws.Range("E24").Value = 5600
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -5200
$ws.Range("E26").Value = 32900
$ws.Range("D27").Value = -5200
$ws.Range("E27").Value = 32900
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 900
$ws.Range("E29").Value = 200
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 22100
$ws.Range("E32").Value = -17400
$ws.Range("D33").Value = -4300
$ws.Range("E33").Value = 33100
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -4300
$ws.Range("E35").Value = 33100
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 59300
$ws.Range("E41").Value = 43500
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 147800
$ws.Range("E43").Value = 93700
$ws.Range("D44").Value = 33500
$ws.Range("E44").Value = 24400
$ws.Range("D45").Value = 18600
$ws.Range("E45").Value = 15700
$ws.Range("D46").Value = 259200
$ws.Range("E46").Value = 177200
$ws.Range("D47").Value = "NA"
$ws.Range("E47").Value = "NA"
$ws.Range("D48").Value = 1555000
$ws.Range("E48").Value = 1226500
$ws.Range("D49").Value = 535400
$ws.Range("E49").Value = 43700
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 121000
$ws.Range("E52").Value = 105600
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 2470600
$ws.Range("E54").Value = 1553100
$ws.Range("D57").Value = 109800
$ws.Range("E57").Value = 101200
$ws.Range("D58").Value = 29700
$ws.Range("E58").Value = 14900
$ws.Range("D59").Value = 89600
$ws.Range("E59").Value = 58000
$ws.Range("D60").Value = 229100
$ws.Range("E60").Value = 174100
$ws.Range("D61").Value = 1575400
$ws.Range("E61").Value = 713300
$ws.Range("D62").Value = 229600
$ws.Range("E62").Value = 195000
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 2034100
$ws.Range("E66").Value = 1082400
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 56100
$ws.Range("E72").Value = 60400
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 436400
$ws.Range("E76").Value = 470700
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -4300
$ws.Range("E81").Value = 33100
$ws.Range("D83").Value = 61100
$ws.Range("E83").Value = 49200
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 97200
$ws.Range("E89").Value = 58800
$ws.Range("D91").Value = -78900
$ws.Range("E91").Value = -63200
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -937100
$ws.Range("E94").Value = -66300
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 855700
$ws.Range("E100").Value = 19300
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 15900
$ws.Range("E102").Value = 11800

# A few quarters further back were also restated; update the cells that now
# live in columns H:I (previously F:G) to the corrected figures.
$ws.Range("H9").Value = 133000
$ws.Range("I9").Value = 99900
$ws.Range("H10").Value = 190000
$ws.Range("I10").Value = 154200
$ws.Range("H17").Value = 289100
$ws.Range("I17").Value = 229600
$ws.Range("H18").Value = 33900
$ws.Range("I18").Value = 24500
$ws.Range("H20").Value = 18000
$ws.Range("I20").Value = -40900
$ws.Range("H32").Value = -18000
$ws.Range("I32").Value = 40900
